# Refresh the cryptocurrency snapshot data (rows 2-51 of the sheet).
# Updated prices/percent-changes, and a handful of rows were re-ranked
# (their Coin/Link/Price/Volume values swapped with the adjacent row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.910.30'
$ws.Range('E2').Value = '  +1.13%  '
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.911.34'
$ws.Range('E3').Value = '  +2.93%  '
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.16%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '471.68'
$ws.Range('E5').Value = '  +10.17%  '
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.38'
$ws.Range('E6').Value = '  +4.97%  '
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('E7').Value = '  +2.22%  '
# Row 8
$ws.Range('E8').Value = '  -0.21%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.744'
$ws.Range('E9').Value = '  -0.18%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.166'
$ws.Range('E10').Value = '  +3.37%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000339'
$ws.Range('E11').Value = '  +2.59%  '
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.32'
$ws.Range('E12').Value = '  -1.02%  '
# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.43'
$ws.Range('E13').Value = '  -3.27%  '
# Row 14
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.521.77'
$ws.Range('E14').Value = '  +2.31%  '
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.04'
$ws.Range('E15').Value = '  -1.16%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.914.62'
$ws.Range('E16').Value = '  +2.93%  '
# Row 17
$ws.Range('E17').Value = '  -0.41%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.07'
$ws.Range('E18').Value = '  -0.62%  '
# Row 19
$ws.Range('E19').Value = '  +2.88%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '67.211.52'
$ws.Range('E20').Value = '  +1.26%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '432.75'
$ws.Range('E21').Value = '  +4.77%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.79'
$ws.Range('E22').Value = '  -3.45%  '
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.35'
$ws.Range('E23').Value = '  +2.13%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.30'
$ws.Range('E24').Value = '  +2.57%  '
# Row 25
$ws.Range('B25').Value = 'EthereumClassic'
$ws.Range('C25').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '38.83'
$ws.Range('E25').Value = '  +4.42%  '
# Row 26
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.56'
$ws.Range('E26').Value = '  +7.00%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.09'
$ws.Range('E27').Value = '  +2.39%  '
# Row 28
$ws.Range('E28').Value = '  +4.44%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.69'
$ws.Range('E29').Value = '  -2.60%  '
# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '736.17'
$ws.Range('E30').Value = '  +4.81%  '
# Row 31
$ws.Range('E31').Value = '  -2.88%  '
# Row 32
$ws.Range('E32').Value = '  +2.89%  '
# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.78'
$ws.Range('E33').Value = '  -0.30%  '
# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '44.22'
$ws.Range('E34').Value = '  +8.98%  '
# Row 35
$ws.Range('E35').Value = '  +3.70%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.12'
$ws.Range('E36').Value = '  +3.42%  '
# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  +0.09%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.38'
$ws.Range('E38').Value = '  -8.48%  '
# Row 39
$ws.Range('B39').Value = 'ThetaToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.22'
$ws.Range('E39').Value = '  +10.89%  '
# Row 40
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0478'
$ws.Range('E40').Value = '  -0.06%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0₃0749'
$ws.Range('E41').Value = '  +8.88%  '
# Row 42
$ws.Range('E42').Value = '  +0.74%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.337'
$ws.Range('E43').Value = '  +3.41%  '
# Row 44
$ws.Range('E44').Value = '  -0.09%  '
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.80'
$ws.Range('E45').Value = '  +6.58%  '
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.19'
$ws.Range('E46').Value = '  +5.62%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.44'
$ws.Range('E47').Value = '  +0.86%  '
# Row 48
$ws.Range('B48').Value = 'Fetch.AI'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.47'
$ws.Range('E48').Value = '  -7.03%  '
# Row 49
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.18'
$ws.Range('E49').Value = '  -0.64%  '
# Row 50
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '26.71'
$ws.Range('E50').Value = '  +2.51%  '
# Row 51
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.89'
$ws.Range('E51').Value = '  +1.73%  '
